$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Remove the two bulleted paragraphs that used to describe
#    "The FeatureReader and FeatureWriter Transformers" / "Integration
#    Transformers" (chapter 3) - the content moved to chapter 4.
# ---------------------------------------------------------------------------
$pFeature = Get-ParaByText $d "The FeatureReader and FeatureWriter Transformers"
$pIntegration = Get-ParaByText $d "Integration Transformers"
$delRange = $d.Range($pFeature.Range.Start, $pIntegration.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 2) The "_GoBack" bookmark used to sit in the "Learn More" paragraph,
#    splitting a " " run and a "t" run apart. Remove that bookmark and
#    restore the paragraph's run layout to a single " t" run followed by
#    the remaining (already separate) runs, matching how Word re-lays the
#    paragraph out once the bookmark that divided it is gone.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$findHost = $d.Content
$findHost.Find.Execute("www.safe.com/training", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterLink = $findHost.End

$findEnd = $d.Content
$findEnd.Find.Execute(" region. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraTailEnd = $findEnd.End

$mergeRange = $d.Range($afterLink, $paraTailEnd)
$mergeXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r w:rsidR="00865029"><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> t</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>o l</w:t></w:r>' + `
    '<w:r w:rsidRPr="003047E1"><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>ocate and register for an upcoming training class in you</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>r</w:t></w:r>' + `
    '<w:r w:rsidRPr="003047E1"><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> region.</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mergeRange.InsertXML($mergeXml)

# ---------------------------------------------------------------------------
# 3) Insert the new "Workflow Transformers" bullet (with a fresh "_GoBack"
#    bookmark, matching Word's behaviour of stamping the last edit point)
#    right before "Managing Attributes" in chapter 4.
# ---------------------------------------------------------------------------
$pManaging = Get-ParaByText $d "Managing Attributes"
$pManaging.Range.InsertParagraphBefore()

$pNew = Get-ParaByText $d "Managing Attributes"
$newPara = $pNew.Previous
$newRange = $newPara.Range
$newRange.Collapse(1)
$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="bulleted-TechBrief"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Workflow </w:t></w:r><w:r><w:t>Transformers</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($newXml)
